$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 37
$ws.Range("A3").Value = 38
$ws.Range("A4:H5").Delete()
